$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

# Paragraph 2 (index 2): (Ref-DJ49SK) -> (Ref-u188936)
Replace-InParagraph 2 "Ref-DJ49SK" "Ref-u188936"

# Paragraph 3 (index 3): (Ref-J6X81N) -> (Ref-u125302), (Ref-29H4Y7) -> (Ref-u125302)
Replace-InParagraph 3 "Ref-J6X81N" "Ref-u125302"
Replace-InParagraph 3 "Ref-29H4Y7" "Ref-u125302"

# Paragraph 4 (index 4): (Ref-SG29DH) -> (Ref-f868999)
Replace-InParagraph 4 "Ref-SG29DH" "Ref-f868999"

# Paragraph 5 (index 5): (Ref-G7H8I9) -> (Ref-s648149), (Ref-J0K1L2) -> (Ref-s648149)
Replace-InParagraph 5 "Ref-G7H8I9" "Ref-s648149"
Replace-InParagraph 5 "Ref-J0K1L2" "Ref-s648149"

# Paragraph 6 (index 6): (Ref-A1B2C3) -> (Ref-f915181), (Ref-D4E5F6) -> (Ref-f915181)
Replace-InParagraph 6 "Ref-A1B2C3" "Ref-f915181"
Replace-InParagraph 6 "Ref-D4E5F6" "Ref-f915181"

# Paragraph 7 (index 7): (Ref-DJ79X2) -> (Al-Sayed, 1998)
Replace-InParagraph 7 "Ref-DJ79X2" "Al-Sayed, 1998"

# Paragraph 8 (index 8): (Ref-A1B2C3) -> (Ref-f155478)
Replace-InParagraph 8 "Ref-A1B2C3" "Ref-f155478"

# Paragraph 9 (index 9): (Ref-J7Y3H2) -> (Ref-s869097)
Replace-InParagraph 9 "Ref-J7Y3H2" "Ref-s869097"
